$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 1983.4736
$ws.Range("J17").Value = 2096.5881
$ws.Range("L17").Value = 6289.7643
$ws.Range("N17").Value = -6625.7643

# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 42593.418
$ws.Range("I28").Value = 53731.05
$ws.Range("K28").Value = 53731.05
$ws.Range("M28").Value = -53246.05

# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 959.8182
$ws.Range("I33").Value = 1164.5555
$ws.Range("J33").Value = 38.5
$ws.Range("K33").Value = 1164.5555
$ws.Range("L33").Value = 38.5
$ws.Range("M33").Value = -935.5554999999999
$ws.Range("N33").Value = -496.5

# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 21742288
$ws.Range("I41").Value = 613.1
$ws.Range("J41").Value = 38466656
$ws.Range("K41").Value = 613.1
$ws.Range("L41").Value = 38466656
$ws.Range("M41").Value = -173.1
$ws.Range("N41").Value = -38467536

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 1889.8
$ws.Range("I43").Value = 1199.8
$ws.Range("K43").Value = 1199.8
$ws.Range("M43").Value = -1130.8

# Row 53 (Leve Item ID 5479)
$ws.Range("H53").Value = 18519578
$ws.Range("I53").Value = 41667444
$ws.Range("J53").Value = 1283.5
$ws.Range("K53").Value = 41667444
$ws.Range("L53").Value = 1283.5
$ws.Range("M53").Value = -41666807
$ws.Range("N53").Value = -2557.5

# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 10422213
$ws.Range("I62").Value = 20836446
$ws.Range("K62").Value = 20836446
$ws.Range("M62").Value = -20835822

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 10422213
$ws.Range("I65").Value = 20836446
$ws.Range("K65").Value = 104182230
$ws.Range("M65").Value = -104179110

# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 545
$ws.Range("I96").Value = 469.16666
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 1407.49998
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = -34.49998000000005
$ws.Range("N96").Value = -5746

# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 45866.773
$ws.Range("I107").Value = 45866.773
$ws.Range("K107").Value = 45866.773
$ws.Range("M107").Value = -43946.773

# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 5851162.5
$ws.Range("I125").Value = 1979
$ws.Range("J125").Value = 9263186
$ws.Range("K125").Value = 17811
$ws.Range("L125").Value = 83368674
$ws.Range("M125").Value = -15351
$ws.Range("N125").Value = -83373594

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 3170
$ws.Range("I132").Value = 3199.9546
$ws.Range("K132").Value = 9599.863799999999
$ws.Range("M132").Value = -7069.863799999999

# Row 134 (Leve Item ID 41997)
$ws.Range("H134").Value = 49998.855
$ws.Range("J134").Value = 49998.855
$ws.Range("L134").Value = 49998.855
$ws.Range("N134").Value = -60138.855

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1114522.8
$ws.Range("I137").Value = 1113595.1
$ws.Range("J137").Value = 1115450.4
$ws.Range("K137").Value = 3340785.3
$ws.Range("L137").Value = 3346351.2
$ws.Range("M137").Value = -3338235.3
$ws.Range("N137").Value = -3351451.2

$ws = $wb.Worksheets.Item("ARM")
# Row 62 (Leve Item ID 10719)
$ws.Range("H62").Value = 29497.5
$ws.Range("J62").Value = 29995
$ws.Range("L62").Value = 29995
$ws.Range("N62").Value = -31243

# Row 65 (Leve Item ID 10719)
$ws.Range("H65").Value = 29497.5
$ws.Range("J65").Value = 29995
$ws.Range("L65").Value = 89985
$ws.Range("N65").Value = -96225

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 886.61536
$ws.Range("I94").Value = 886.61536
$ws.Range("K94").Value = 886.61536
$ws.Range("M94").Value = -435.61536

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 33054.09
$ws.Range("I134").Value = 1955
$ws.Range("K134").Value = 5865
$ws.Range("M134").Value = -3330

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 31994.266
$ws.Range("I31").Value = 1730.8846
$ws.Range("J31").Value = 130350.25
$ws.Range("K31").Value = 1730.8846
$ws.Range("L31").Value = 130350.25
$ws.Range("M31").Value = -1435.8846
$ws.Range("N31").Value = -130940.25

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 31994.266
$ws.Range("I34").Value = 1730.8846
$ws.Range("J34").Value = 130350.25
$ws.Range("K34").Value = 1730.8846
$ws.Range("L34").Value = 130350.25
$ws.Range("M34").Value = -1528.8846
$ws.Range("N34").Value = -130754.25

# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 6373.2354
$ws.Range("I99").Value = 5444.6
$ws.Range("J99").Value = 6760.1665
$ws.Range("K99").Value = 5444.6
$ws.Range("L99").Value = 6760.1665
$ws.Range("M99").Value = -3946.6
$ws.Range("N99").Value = -9756.166499999999

# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 4550.7646
$ws.Range("I122").Value = 3038.1428
$ws.Range("K122").Value = 9114.428400000001
$ws.Range("M122").Value = -6664.428400000001

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 6373.2354
$ws.Range("I126").Value = 5444.6
$ws.Range("J126").Value = 6760.1665
$ws.Range("K126").Value = 16333.8
$ws.Range("L126").Value = 20280.4995
$ws.Range("M126").Value = -13863.8
$ws.Range("N126").Value = -25220.4995

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 3919.5862
$ws.Range("I132").Value = 3280.55
$ws.Range("K132").Value = 9841.650000000001
$ws.Range("M132").Value = -7311.650000000001

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 669903.3
$ws.Range("I134").Value = 457553
$ws.Range("K134").Value = 1372659
$ws.Range("M134").Value = -1370124

$ws = $wb.Worksheets.Item("GSM")
# Row 98 (Leve Item ID 18359)
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990

# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 4278.2856
$ws.Range("I102").Value = 505
$ws.Range("J102").Value = 4907.1665
$ws.Range("K102").Value = 505
$ws.Range("L102").Value = 4907.1665
$ws.Range("M102").Value = 1117
$ws.Range("N102").Value = -8151.1665

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 3975.182
$ws.Range("I126").Value = 3145.8
$ws.Range("K126").Value = 9437.400000000001
$ws.Range("M126").Value = -6967.400000000001

# Row 130 (Leve Item ID 34692)
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 505715.78
$ws.Range("I132").Value = 1006423.4
$ws.Range("J132").Value = 88459.414
$ws.Range("K132").Value = 3019270.2
$ws.Range("L132").Value = 265378.242
$ws.Range("M132").Value = -3016740.2
$ws.Range("N132").Value = -270438.242

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 230739.4
$ws.Range("I7").Value = 3431.923
$ws.Range("J7").Value = 559072.4399999999
$ws.Range("K7").Value = 3431.923
$ws.Range("L7").Value = 559072.4399999999
$ws.Range("M7").Value = -3319.923
$ws.Range("N7").Value = -559296.4399999999

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 1455.7273
$ws.Range("I22").Value = 1074
$ws.Range("J22").Value = 1913.8
$ws.Range("K22").Value = 1074
$ws.Range("L22").Value = 1913.8
$ws.Range("M22").Value = -779
$ws.Range("N22").Value = -2503.8

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 1455.7273
$ws.Range("I27").Value = 1074
$ws.Range("J27").Value = 1913.8
$ws.Range("K27").Value = 1074
$ws.Range("L27").Value = 1913.8
$ws.Range("M27").Value = -967
$ws.Range("N27").Value = -2127.8

# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 2003350
$ws.Range("I40").Value = 2943886.5
$ws.Range("J40").Value = 4709.875
$ws.Range("K40").Value = 2943886.5
$ws.Range("L40").Value = 4709.875
$ws.Range("M40").Value = -2943750.5
$ws.Range("N40").Value = -4981.875

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 624.45715
$ws.Range("I55").Value = 180.2963
$ws.Range("J55").Value = 2123.5
$ws.Range("K55").Value = 180.2963
$ws.Range("L55").Value = 2123.5
$ws.Range("M55").Value = -7.296300000000002
$ws.Range("N55").Value = -2469.5

# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 2007.7273
$ws.Range("I93").Value = 1837.4445
$ws.Range("J93").Value = 2774
$ws.Range("K93").Value = 1837.4445
$ws.Range("L93").Value = 2774
$ws.Range("M93").Value = -589.4445000000001
$ws.Range("N93").Value = -5270

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 230739.4
$ws.Range("I126").Value = 3431.923
$ws.Range("J126").Value = 559072.4399999999
$ws.Range("K126").Value = 10295.769
$ws.Range("L126").Value = 1677217.32
$ws.Range("M126").Value = -7825.769
$ws.Range("N126").Value = -1682157.32

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3553.9048
$ws.Range("I132").Value = 2476.2144
$ws.Range("K132").Value = 7428.6432
$ws.Range("M132").Value = -4898.6432

# Row 138 (Leve Item ID 42334)
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 1921.5625
$ws.Range("I81").Value = 1403.5834
$ws.Range("J81").Value = 3475.5
$ws.Range("K81").Value = 2807.1668
$ws.Range("L81").Value = 6951
$ws.Range("M81").Value = -1746.1668
$ws.Range("N81").Value = -9073

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 1921.5625
$ws.Range("I84").Value = 1403.5834
$ws.Range("J84").Value = 3475.5
$ws.Range("K84").Value = 14035.834
$ws.Range("L84").Value = 34755
$ws.Range("M84").Value = -8731.833999999999
$ws.Range("N84").Value = -45363

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 4671.75
$ws.Range("I126").Value = 3200
$ws.Range("K126").Value = 9600
$ws.Range("M126").Value = -7130

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 43538.6
$ws.Range("I132").Value = 2516.611
$ws.Range("K132").Value = 7549.833
$ws.Range("M132").Value = -5019.833

# Row 138 (Leve Item ID 42347)
$ws.Range("H138").Value = 80770
$ws.Range("J138").Value = 80770
$ws.Range("L138").Value = 80770
$ws.Range("N138").Value = -91050
